# Remove the placeholder "-" values from the "Not EN 388 rated" / incomplete rows
# so the cells are genuinely empty instead of containing a literal dash.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J3").ClearContents()
$ws.Range("F5:J5").ClearContents()
$ws.Range("F6:J6").ClearContents()
$ws.Range("F7:J7").ClearContents()
$ws.Range("F8:J8").ClearContents()
$ws.Range("F14:J14").ClearContents()
$ws.Range("J18").ClearContents()
$ws.Range("J19").ClearContents()

# Clear the (redundant, identical-to-default) explicit formatting on the
# Image column header + blank cells so the duplicate style gets dropped.
$ws.Range("D1:D16").ClearFormats()
